$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD675H57H%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T010004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJHMEUCIQCAxpDCI0%2F8aKBWKU23Ug3VxeZqNQFvWPzf%2B6s32P1w7QIgbHml%2BO4uZGZkTMxn1hupIeKDN1DzlBmnH7GOb0NbHUMq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGf9I%2B70MozV5rw6iCq6AzUR26dABzsrrJEqhN20YkwDSdTLVooattHNUqLHjyDoqLzz534kKt7dZMKbcYK5zXx7EzfBfCskfVDh8YAogZKp9CUxTfFKuDa%2By%2BlR3cMITQnsupiJHzobLAcrUHmDbKfDlK0ZQVyOF6GMdzG%2FAzdWrGlpNqVEy1xfTdBDxSsQhJ2ASTa0KEmjzoJpy%2BMm4vjZFnn%2FdUtGiK89AmRO8caIsb2V%2FkG0OJekW0mG%2Bjd9atreGFwEqh6qCdIXquTxBcnxeeWWgQG3ZVdVeKak0dRtRfRjsa9CBwPBmhLgxRbNHxunIeLlvSdin1OkJ0JmQecg7aHuots3TGFWcLVpNP0tsdqGAzv7R9ZmpoHGoEFjnZp%2FsLZvQIviuDaiSmqd7F0crnbeaL7fgGGutA0%2F24ULbs9Um9FKZjI5F4pDZcWDAh8M4bLTA9x0hR7kmAN0wPYnKpH6K9eYXqEUck1ODNkTGxcPrSTtKisMgcnT7SviT1Q9UrZFaZ1oPSfCU%2B6Fnezf2NKoBy73cJwAmewEriqq%2BJSySHb78s3lRLeitCDcV7C87dBzaaNGi8xQ0l09c43bBo%2Fl8wqkXosw4uymvwY6pQECfUR0E8DjHWZ4NXhhcdo0scTIfBDRdnNM3lBhCnidUdYzbOmYYvLyBZgFXesQnyjwuRSPKlN2PxRQrGWKPKT3Yj3b00MIWaEZuBE5pDnYh9ZD0hhJYeIfVWr77bfO5QpyBvzKRALdILMGIQ35QCALmJ6uG4ibdTSUuWDCFwBkIXmEUJT5bWqKDrbTHiaY7bibfkCSz%2BDjTvcACMmT41TTZh%2FXI3A%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2009%3A00%3A04%20GMT&X-Amz-Signature=8127ae4d7e78928bdc75eb6e86fe74d1b284d84c701dc075a5d0dc6e75868857'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD675H57H%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T010004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJHMEUCIQCAxpDCI0%2F8aKBWKU23Ug3VxeZqNQFvWPzf%2B6s32P1w7QIgbHml%2BO4uZGZkTMxn1hupIeKDN1DzlBmnH7GOb0NbHUMq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGf9I%2B70MozV5rw6iCq6AzUR26dABzsrrJEqhN20YkwDSdTLVooattHNUqLHjyDoqLzz534kKt7dZMKbcYK5zXx7EzfBfCskfVDh8YAogZKp9CUxTfFKuDa%2By%2BlR3cMITQnsupiJHzobLAcrUHmDbKfDlK0ZQVyOF6GMdzG%2FAzdWrGlpNqVEy1xfTdBDxSsQhJ2ASTa0KEmjzoJpy%2BMm4vjZFnn%2FdUtGiK89AmRO8caIsb2V%2FkG0OJekW0mG%2Bjd9atreGFwEqh6qCdIXquTxBcnxeeWWgQG3ZVdVeKak0dRtRfRjsa9CBwPBmhLgxRbNHxunIeLlvSdin1OkJ0JmQecg7aHuots3TGFWcLVpNP0tsdqGAzv7R9ZmpoHGoEFjnZp%2FsLZvQIviuDaiSmqd7F0crnbeaL7fgGGutA0%2F24ULbs9Um9FKZjI5F4pDZcWDAh8M4bLTA9x0hR7kmAN0wPYnKpH6K9eYXqEUck1ODNkTGxcPrSTtKisMgcnT7SviT1Q9UrZFaZ1oPSfCU%2B6Fnezf2NKoBy73cJwAmewEriqq%2BJSySHb78s3lRLeitCDcV7C87dBzaaNGi8xQ0l09c43bBo%2Fl8wqkXosw4uymvwY6pQECfUR0E8DjHWZ4NXhhcdo0scTIfBDRdnNM3lBhCnidUdYzbOmYYvLyBZgFXesQnyjwuRSPKlN2PxRQrGWKPKT3Yj3b00MIWaEZuBE5pDnYh9ZD0hhJYeIfVWr77bfO5QpyBvzKRALdILMGIQ35QCALmJ6uG4ibdTSUuWDCFwBkIXmEUJT5bWqKDrbTHiaY7bibfkCSz%2BDjTvcACMmT41TTZh%2FXI3A%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2009%3A00%3A04%20GMT&X-Amz-Signature=4e1812c528aae657857880b4de202c9ac8d26ea8e861674f442e3978051ba5f8'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD675H57H%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T010004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJHMEUCIQCAxpDCI0%2F8aKBWKU23Ug3VxeZqNQFvWPzf%2B6s32P1w7QIgbHml%2BO4uZGZkTMxn1hupIeKDN1DzlBmnH7GOb0NbHUMq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGf9I%2B70MozV5rw6iCq6AzUR26dABzsrrJEqhN20YkwDSdTLVooattHNUqLHjyDoqLzz534kKt7dZMKbcYK5zXx7EzfBfCskfVDh8YAogZKp9CUxTfFKuDa%2By%2BlR3cMITQnsupiJHzobLAcrUHmDbKfDlK0ZQVyOF6GMdzG%2FAzdWrGlpNqVEy1xfTdBDxSsQhJ2ASTa0KEmjzoJpy%2BMm4vjZFnn%2FdUtGiK89AmRO8caIsb2V%2FkG0OJekW0mG%2Bjd9atreGFwEqh6qCdIXquTxBcnxeeWWgQG3ZVdVeKak0dRtRfRjsa9CBwPBmhLgxRbNHxunIeLlvSdin1OkJ0JmQecg7aHuots3TGFWcLVpNP0tsdqGAzv7R9ZmpoHGoEFjnZp%2FsLZvQIviuDaiSmqd7F0crnbeaL7fgGGutA0%2F24ULbs9Um9FKZjI5F4pDZcWDAh8M4bLTA9x0hR7kmAN0wPYnKpH6K9eYXqEUck1ODNkTGxcPrSTtKisMgcnT7SviT1Q9UrZFaZ1oPSfCU%2B6Fnezf2NKoBy73cJwAmewEriqq%2BJSySHb78s3lRLeitCDcV7C87dBzaaNGi8xQ0l09c43bBo%2Fl8wqkXosw4uymvwY6pQECfUR0E8DjHWZ4NXhhcdo0scTIfBDRdnNM3lBhCnidUdYzbOmYYvLyBZgFXesQnyjwuRSPKlN2PxRQrGWKPKT3Yj3b00MIWaEZuBE5pDnYh9ZD0hhJYeIfVWr77bfO5QpyBvzKRALdILMGIQ35QCALmJ6uG4ibdTSUuWDCFwBkIXmEUJT5bWqKDrbTHiaY7bibfkCSz%2BDjTvcACMmT41TTZh%2FXI3A%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2009%3A00%3A04%20GMT&X-Amz-Signature=f35b21ed78de389a01d52f15ed710c1a434a133e817866f5623a01619844ed0c'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD675H57H%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T010005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJHMEUCIQCAxpDCI0%2F8aKBWKU23Ug3VxeZqNQFvWPzf%2B6s32P1w7QIgbHml%2BO4uZGZkTMxn1hupIeKDN1DzlBmnH7GOb0NbHUMq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGf9I%2B70MozV5rw6iCq6AzUR26dABzsrrJEqhN20YkwDSdTLVooattHNUqLHjyDoqLzz534kKt7dZMKbcYK5zXx7EzfBfCskfVDh8YAogZKp9CUxTfFKuDa%2By%2BlR3cMITQnsupiJHzobLAcrUHmDbKfDlK0ZQVyOF6GMdzG%2FAzdWrGlpNqVEy1xfTdBDxSsQhJ2ASTa0KEmjzoJpy%2BMm4vjZFnn%2FdUtGiK89AmRO8caIsb2V%2FkG0OJekW0mG%2Bjd9atreGFwEqh6qCdIXquTxBcnxeeWWgQG3ZVdVeKak0dRtRfRjsa9CBwPBmhLgxRbNHxunIeLlvSdin1OkJ0JmQecg7aHuots3TGFWcLVpNP0tsdqGAzv7R9ZmpoHGoEFjnZp%2FsLZvQIviuDaiSmqd7F0crnbeaL7fgGGutA0%2F24ULbs9Um9FKZjI5F4pDZcWDAh8M4bLTA9x0hR7kmAN0wPYnKpH6K9eYXqEUck1ODNkTGxcPrSTtKisMgcnT7SviT1Q9UrZFaZ1oPSfCU%2B6Fnezf2NKoBy73cJwAmewEriqq%2BJSySHb78s3lRLeitCDcV7C87dBzaaNGi8xQ0l09c43bBo%2Fl8wqkXosw4uymvwY6pQECfUR0E8DjHWZ4NXhhcdo0scTIfBDRdnNM3lBhCnidUdYzbOmYYvLyBZgFXesQnyjwuRSPKlN2PxRQrGWKPKT3Yj3b00MIWaEZuBE5pDnYh9ZD0hhJYeIfVWr77bfO5QpyBvzKRALdILMGIQ35QCALmJ6uG4ibdTSUuWDCFwBkIXmEUJT5bWqKDrbTHiaY7bibfkCSz%2BDjTvcACMmT41TTZh%2FXI3A%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2009%3A00%3A05%20GMT&X-Amz-Signature=847e9efd1dda10ae0dae5354897117eac408857c96cf073792ccb4105ce9a487'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD675H57H%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T010004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJHMEUCIQCAxpDCI0%2F8aKBWKU23Ug3VxeZqNQFvWPzf%2B6s32P1w7QIgbHml%2BO4uZGZkTMxn1hupIeKDN1DzlBmnH7GOb0NbHUMq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGf9I%2B70MozV5rw6iCq6AzUR26dABzsrrJEqhN20YkwDSdTLVooattHNUqLHjyDoqLzz534kKt7dZMKbcYK5zXx7EzfBfCskfVDh8YAogZKp9CUxTfFKuDa%2By%2BlR3cMITQnsupiJHzobLAcrUHmDbKfDlK0ZQVyOF6GMdzG%2FAzdWrGlpNqVEy1xfTdBDxSsQhJ2ASTa0KEmjzoJpy%2BMm4vjZFnn%2FdUtGiK89AmRO8caIsb2V%2FkG0OJekW0mG%2Bjd9atreGFwEqh6qCdIXquTxBcnxeeWWgQG3ZVdVeKak0dRtRfRjsa9CBwPBmhLgxRbNHxunIeLlvSdin1OkJ0JmQecg7aHuots3TGFWcLVpNP0tsdqGAzv7R9ZmpoHGoEFjnZp%2FsLZvQIviuDaiSmqd7F0crnbeaL7fgGGutA0%2F24ULbs9Um9FKZjI5F4pDZcWDAh8M4bLTA9x0hR7kmAN0wPYnKpH6K9eYXqEUck1ODNkTGxcPrSTtKisMgcnT7SviT1Q9UrZFaZ1oPSfCU%2B6Fnezf2NKoBy73cJwAmewEriqq%2BJSySHb78s3lRLeitCDcV7C87dBzaaNGi8xQ0l09c43bBo%2Fl8wqkXosw4uymvwY6pQECfUR0E8DjHWZ4NXhhcdo0scTIfBDRdnNM3lBhCnidUdYzbOmYYvLyBZgFXesQnyjwuRSPKlN2PxRQrGWKPKT3Yj3b00MIWaEZuBE5pDnYh9ZD0hhJYeIfVWr77bfO5QpyBvzKRALdILMGIQ35QCALmJ6uG4ibdTSUuWDCFwBkIXmEUJT5bWqKDrbTHiaY7bibfkCSz%2BDjTvcACMmT41TTZh%2FXI3A%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2009%3A00%3A04%20GMT&X-Amz-Signature=d06543b461602e96662f4caba01a2cacc389d4cb2d5fe95c66d437ff8cd4fbe3'
